$wb = $excel.ActiveWorkbook

$newUuid   = "d6745a91-71b4-416f-93c0-c571990ec9a6"
$xlfCommit = "12766b8e48f7afdf62a9396375eed44b8931d39c"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/e7f6d9994f2051e4e5b6632d66b08c11d6e148ac/e2e/$newUuid.md"
$zhXlfName = "$newUuid.$xlfCommit.zh-cn.xlf"
$deXlfName = "$newUuid.$xlfCommit.de-de.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5470b9148651a35570e00f6587938f693abf307a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce25f650c1c002d2cafda1d3e373441db12b3a77/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

# --- Sheet 1: Overview ---
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrl, "", "", "$newUuid.md")
$wsOverview.Range("B4").Value2 = "Ready for handoff"
$wsOverview.Range("C4").Value2 = "Ready for handoff"
$wsOverview.Range("D4").Value2 = "2016-50-13 16:50:38"

# --- Sheet 2: zh-cn ---
$wsZh = $wb.Worksheets.Item(2)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdUrl, "", "", "$newUuid.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), $mdUrl, "", "", ".md")
$wsZh.Range("C4").Value2 = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $zhXlfUrl, "", "", $zhXlfName)
$wsZh.Range("E4").Value2 = "2016-03-13 16:50:35"
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value2 = "Include"

# --- Sheet 3: de-de ---
$wsDe = $wb.Worksheets.Item(3)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdUrl, "", "", "$newUuid.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), $mdUrl, "", "", ".md")
$wsDe.Range("C4").Value2 = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $deXlfUrl, "", "", $deXlfName)
$wsDe.Range("E4").Value2 = "2016-03-13 16:50:38"
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value2 = "Include"

Write-Host "Done"
